$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 866.2143
$ws.Range("I38").Value = 191.8
$ws.Range("J38").Value = 2552.25
$ws.Range("K38").Value = 575.4000000000001
$ws.Range("L38").Value = 7656.75
$ws.Range("M38").Value = -203.4000000000001
$ws.Range("N38").Value = -8400.75

$ws.Range("H39").Value = 66.27273
$ws.Range("I39").Value = 66.27273
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 198.81819
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 97.18181000000001
$ws.Range("N39").Value = $null

$ws.Range("H42").Value = 294.5
$ws.Range("J42").Value = 436.33334
$ws.Range("L42").Value = 1309.00002
$ws.Range("N42").Value = -1769.00002

$ws.Range("H98").Value = 339938.12
$ws.Range("I98").Value = 430471.47
$ws.Range("K98").Value = 430471.47
$ws.Range("M98").Value = -428973.47

$ws.Range("H116").Value = 5537695.5
$ws.Range("I116").Value = 9884418
$ws.Range("J116").Value = 5503
$ws.Range("K116").Value = 9884418
$ws.Range("L116").Value = 5503
$ws.Range("M116").Value = -9880976
$ws.Range("N116").Value = -12387

$ws.Range("H122").Value = 339938.12
$ws.Range("I122").Value = 430471.47
$ws.Range("K122").Value = 1291414.41
$ws.Range("M122").Value = -1288964.41

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2724.75
$ws.Range("I32").Value = 1709.2858
$ws.Range("K32").Value = 1709.2858
$ws.Range("M32").Value = -1422.2858

$ws.Range("H45").Value = 1682.3529
$ws.Range("I45").Value = 1373.3334
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1373.3334
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -996.3334
$ws.Range("N45").Value = -4754

$ws.Range("H132").Value = 1312.3704
$ws.Range("I132").Value = 957
$ws.Range("K132").Value = 2871
$ws.Range("M132").Value = -341

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3027.6287
$ws.Range("I134").Value = 1308.55
$ws.Range("J134").Value = 5319.7334
$ws.Range("K134").Value = 3925.65
$ws.Range("L134").Value = 15959.2002
$ws.Range("M134").Value = -1390.65
$ws.Range("N134").Value = -21029.2002

$ws.Range("H137").Value = 65926.664
$ws.Range("J137").Value = 65926.664
$ws.Range("L137").Value = 65926.664
$ws.Range("N137").Value = -76126.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1356.0822
$ws.Range("I31").Value = 919.75555
$ws.Range("J31").Value = 2057.3215
$ws.Range("K31").Value = 919.75555
$ws.Range("L31").Value = 2057.3215
$ws.Range("M31").Value = -624.75555
$ws.Range("N31").Value = -2647.3215

$ws.Range("H34").Value = 1356.0822
$ws.Range("I34").Value = 919.75555
$ws.Range("J34").Value = 2057.3215
$ws.Range("K34").Value = 919.75555
$ws.Range("L34").Value = 2057.3215
$ws.Range("M34").Value = -717.75555
$ws.Range("N34").Value = -2461.3215

$ws.Range("H58").Value = 2147.652
$ws.Range("I58").Value = 1113
$ws.Range("J58").Value = 2512.8235
$ws.Range("K58").Value = 1113
$ws.Range("L58").Value = 2512.8235
$ws.Range("M58").Value = -910
$ws.Range("N58").Value = -2918.8235

$ws.Range("H132").Value = 1977.4166
$ws.Range("I132").Value = 1319.1072
$ws.Range("J132").Value = 4281.5
$ws.Range("K132").Value = 3957.3216
$ws.Range("L132").Value = 12844.5
$ws.Range("M132").Value = -1427.3216
$ws.Range("N132").Value = -17904.5

$ws.Range("H134").Value = 2190.1482
$ws.Range("I134").Value = 798.7222
$ws.Range("J134").Value = 4973
$ws.Range("K134").Value = 2396.1666
$ws.Range("L134").Value = 14919
$ws.Range("M134").Value = 138.8334
$ws.Range("N134").Value = -19989

$ws.Range("H136").Value = 2147.652
$ws.Range("I136").Value = 1113
$ws.Range("J136").Value = 2512.8235
$ws.Range("K136").Value = 3339
$ws.Range("L136").Value = 7538.470499999999
$ws.Range("M136").Value = -789
$ws.Range("N136").Value = -12638.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2419.247
$ws.Range("I68").Value = 2684.127
$ws.Range("J68").Value = 1660.7273
$ws.Range("K68").Value = 8052.380999999999
$ws.Range("L68").Value = 4982.1819
$ws.Range("M68").Value = -7241.380999999999
$ws.Range("N68").Value = -6604.1819

$ws.Range("H71").Value = 2419.247
$ws.Range("I71").Value = 2684.127
$ws.Range("J71").Value = 1660.7273
$ws.Range("K71").Value = 24157.143
$ws.Range("L71").Value = 14946.5457
$ws.Range("M71").Value = -20101.143
$ws.Range("N71").Value = -23058.5457

$ws.Range("H104").Value = 5000
$ws.Range("J104").Value = 5000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -20242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 666.8333
$ws.Range("I22").Value = 711.6923
$ws.Range("J22").Value = 550.2
$ws.Range("K22").Value = 711.6923
$ws.Range("L22").Value = 550.2
$ws.Range("M22").Value = -416.6923
$ws.Range("N22").Value = -1140.2

$ws.Range("H27").Value = 666.8333
$ws.Range("I27").Value = 711.6923
$ws.Range("J27").Value = 550.2
$ws.Range("K27").Value = 711.6923
$ws.Range("L27").Value = 550.2
$ws.Range("M27").Value = -604.6923
$ws.Range("N27").Value = -764.2

$ws.Range("H46").Value = 589.9
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 544.3333
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 544.3333
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -920.3333

$ws.Range("H55").Value = 733.6667
$ws.Range("I55").Value = 899
$ws.Range("J55").Value = 651
$ws.Range("K55").Value = 899
$ws.Range("L55").Value = 651
$ws.Range("M55").Value = -726
$ws.Range("N55").Value = -997

$ws.Range("H93").Value = 761.4706
$ws.Range("I93").Value = 727.5625
$ws.Range("J93").Value = 1304
$ws.Range("K93").Value = 727.5625
$ws.Range("L93").Value = 1304
$ws.Range("M93").Value = 520.4375
$ws.Range("N93").Value = -3800

$ws.Range("H100").Value = 2332.2258
$ws.Range("I100").Value = 1745.3636
$ws.Range("J100").Value = 2655
$ws.Range("K100").Value = 1745.3636
$ws.Range("L100").Value = 2655
$ws.Range("M100").Value = -1204.3636
$ws.Range("N100").Value = -3737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1317687.5
$ws.Range("J81").Value = 5928.5713
$ws.Range("L81").Value = 11857.1426
$ws.Range("N81").Value = -13979.1426

$ws.Range("H84").Value = 1317687.5
$ws.Range("J84").Value = 5928.5713
$ws.Range("L84").Value = 59285.713
$ws.Range("N84").Value = -69893.71299999999

$ws.Range("H107").Value = 5051551.5
$ws.Range("I107").Value = 13890114
$ws.Range("J107").Value = 944.4286
$ws.Range("K107").Value = 41670342
$ws.Range("L107").Value = 2833.2858
$ws.Range("M107").Value = -41668422
$ws.Range("N107").Value = -6673.2858

$ws.Range("H136").Value = 6173616.5
$ws.Range("I136").Value = 8333757
$ws.Range("J136").Value = 1787
$ws.Range("K136").Value = 25001271
$ws.Range("L136").Value = 5361
$ws.Range("M136").Value = -24998721
$ws.Range("N136").Value = -10461
